# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps for the b4de5971-... row across
# the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column G, row 2
# (shares its shared-string value with de-de!H2, so both move together)
$wsOverview.Range("G2").Value = "2016-08-25 09:09:48"

# zh-cn sheet: "Correspond Handoff Datetime" column H, row 2
$wsZhCn.Range("H2").Value = "2016-08-25 09:09:42"
# zh-cn sheet: "Correspond Handback DateTime" column K, row 2
$wsZhCn.Range("K2").Value = "2016-08-25 09:10:25"

# de-de sheet: "Correspond Handoff Datetime" column H, row 2
# (same underlying value as Overview!G2 before the edit)
$wsDeDe.Range("H2").Value = "2016-08-25 09:09:48"
# de-de sheet: "Correspond Handback DateTime" column K, row 2
$wsDeDe.Range("K2").Value = "2016-08-25 09:10:33"
